$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet right after "总计" (i.e. before "2022-Q1")
# ---------------------------------------------------------------------------
$anchorSheet = $wb.Worksheets.Item("2022-Q1")
$q3 = $wb.Worksheets.Add($anchorSheet)
$q3.Name = "2022-Q3"

# Match the page margins used by the rest of the workbook's data sheets
# (0.75in / 1in / 0.5in instead of Excel's factory defaults).
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Header row (bold / bordered / centered) - copy formatting from a sibling sheet
$wb.Worksheets.Item("2021-Q4").Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Column A row-index style (centered / bordered) - copy formatting from a sibling sheet
$wb.Worksheets.Item("2021-Q4").Range("A2").Copy()
$q3.Range("A2:A4").PasteSpecial(-4122)
$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(4, 1).Value = 2

# Data rows. Fund code / scale / position figures are stored as plain TEXT in
# this workbook (even though they look numeric), so they are entered with a
# leading apostrophe to stop them being auto-converted into numbers.
$q3.Cells.Item(2, 2).Value = "'159617"
$q3.Cells.Item(2, 3).Value = "华夏中证智选500价值稳健策略ETF"
$q3.Cells.Item(2, 4).Value = "'2.93"
$q3.Cells.Item(2, 5).Value = "'97.05"
$q3.Cells.Item(2, 6).Value = "'1.41"
$q3.Cells.Item(2, 7).Value = "'0.0413"
$q3.Cells.Item(2, 8).Value = 6

$q3.Cells.Item(3, 2).Value = "'006347"
$q3.Cells.Item(3, 3).Value = "安信量化优选股票C"
$q3.Cells.Item(3, 4).Value = "'0.15"
$q3.Cells.Item(3, 5).Value = "'90.50"
$q3.Cells.Item(3, 6).Value = "'0.77"
$q3.Cells.Item(3, 7).Value = "'0.0012"
$q3.Cells.Item(3, 8).Value = 9

$q3.Cells.Item(4, 2).Value = "'006346"
$q3.Cells.Item(4, 3).Value = "安信量化优选股票A"
$q3.Cells.Item(4, 4).Value = "'0.03"
$q3.Cells.Item(4, 5).Value = "'90.50"
$q3.Cells.Item(4, 6).Value = "'0.77"
$q3.Cells.Item(4, 7).Value = "'0.0002"
$q3.Cells.Item(4, 8).Value = 9

# The apostrophe entry marks these cells with a "quote prefix" style; clear
# that cosmetic flag by re-pasting formats from an untouched, unstyled cell.
$q3.Cells.Item(10, 1).Copy()
$q3.Range("B2:G4").PasteSpecial(-4122)
$q3.Cells.Item(10, 1).ClearContents()

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: a new 2022-Q3 row is inserted at the
#    top of the data table, pushing every existing quarter down by one row
#    and adding 2020-Q4 at the new row 8.
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Shift existing data rows (2..7) down to (3..8), bottom-up so values aren't
# clobbered before they are read.
for ($r = 7; $r -ge 2; $r--) {
    $tot.Cells.Item($r + 1, 2).Value = $tot.Cells.Item($r, 2).Value2
    $tot.Cells.Item($r + 1, 3).Value = $tot.Cells.Item($r, 3).Value2
    $tot.Cells.Item($r + 1, 4).Value = $tot.Cells.Item($r, 4).Value2
}

# New row 2: 2022-Q3 figures
$tot.Cells.Item(2, 2).Value = "2022-Q3"
$tot.Cells.Item(2, 3).Value = 3
$tot.Cells.Item(2, 4).Value = 0.04

# Extend the row-index column A down to the new row 8, matching the style
# already used by the other index cells in column A.
$tot.Cells.Item(7, 1).Copy()
$tot.Cells.Item(8, 1).PasteSpecial(-4122)
$tot.Cells.Item(8, 1).Value = 6
